# Updates attendance summary columns (Invalid = G, Absent = H) on Sheet1.
# Most attendance dates had no marked record, so the person is counted as
# "Absent" (H = 1). Two dates (row 3 and row 6) were additionally marked
# as "Invalid" (G = 1). Row 12 (01/09/2022) instead had an actual
# attendance record, so it is counted towards "Total Attendance Count"
# and "Real" (D = 1, E = 1) rather than Absent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Absent
$ws.Range("H4").Value = 1

# Row 5: Absent
$ws.Range("H5").Value = 1

# Row 6: Invalid + Absent
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

# Row 7: Absent
$ws.Range("H7").Value = 1

# Row 8: Absent
$ws.Range("H8").Value = 1

# Row 9: Absent
$ws.Range("H9").Value = 1

# Row 10: Absent
$ws.Range("H10").Value = 1

# Row 11: Absent
$ws.Range("H11").Value = 1

# Row 12: Total Attendance Count + Real (present, not absent)
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: Absent
$ws.Range("H13").Value = 1

# Row 14: Absent
$ws.Range("H14").Value = 1

# Row 15: Absent
$ws.Range("H15").Value = 1

# Row 16: Absent
$ws.Range("H16").Value = 1

# Row 17: Absent
$ws.Range("H17").Value = 1

# Row 18: Absent
$ws.Range("H18").Value = 1
